$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value pairs derived from the updated TPM recomputation.
$updates = @"
G2=28.09534366666667
H2=84.28603100000001
I2=0.6431807885786103
J2=0.6438169333507339
M2=0.5804443333333333
N2=1.741333
O2=0.002431273010151717
P2=0.002435427107574628
Q2=16.30778302436922
R2=146.770047219323
S2=0.001563748091919273
T2=0.001567969211797945
G3=28.09534366666667
H3=84.28603100000001
I3=0.6431807885786103
J3=0.6438169333507339
O3=0.0004752041289926495
P3=0.00047601606752829
S3=0.0003056421664213039
T3=0.0003064672048217395
G4=28.09534366666667
H4=84.28603100000001
I4=0.6431807885786103
J4=0.6438169333507339
M4=136.1000366666667
N4=408.30011
O4=0.5700742118164518
P4=0.5710482463260632
Q4=3823.777303195935
R4=34413.99572876342
S4=0.3666607811044352
T4=0.3676505307449605
G5=28.09534366666667
H5=84.28603100000001
I5=0.6431807885786103
J5=0.6438169333507339
M5=1.221658
N5=2.443316
O5=0.005117086949542552
P5=0.003417220037046797
Q5=34.32290135313266
R5=205.937408118796
S5=0.003291212019432094
T5=0.00220006412483615
G6=28.09534366666667
H6=84.28603100000001
I6=0.6431807885786103
J6=0.6438169333507339
M6=100.7253213333333
N6=302.175964
O6=0.4219022240948613
P6=0.4226230904617871
Q6=2829.912518795432
R6=25469.21266915889
S6=0.2713594051964024
T6=0.2720919020643175
E7=3
F7=1
G7=15.16980166666667
H7=45.509405
I7=0.347279076358968
J7=0.3476225564081497
M7=0.5804443333333333
N7=1.741333
O7=0.002431273010151717
P7=0.002435427107574628
Q7=8.805225415207222
R7=79.247028736865
S7=0.0008443302453419764
T7=0.0008466093970807981
E8=3
F8=1
G8=15.16980166666667
H8=45.509405
I8=0.347279076358968
J8=0.3476225564081497
O8=0.0004752041289926495
P8=0.00047601606752829
Q8=1.721024112284444
R8=15.48921701056
S8=0.0001650284509985352
T8=0.0001654739222855386
E9=3
F9=1
G9=15.16980166666667
H9=45.509405
I9=0.347279076358968
J9=0.3476225564081497
M9=136.1000366666667
N9=408.30011
O9=0.5700742118164518
P9=0.5710482463260632
Q9=2064.610563059395
R9=18581.49506753455
S9=0.1979748457356841
T9=0.1985092512202568
E10=3
F10=1
G10=15.16980166666667
H10=45.509405
I10=0.347279076358968
J10=0.3476225564081497
M10=1.221658
N10=2.443316
O10=0.005117086949542552
P10=0.003417220037046797
Q10=18.53230956449666
R10=111.19385738698
S10=0.001777057229485667
T10=0.001187902765087359
E11=3
F11=1
G11=15.16980166666667
H11=45.509405
I11=0.347279076358968
J11=0.3476225564081497
M11=100.7253213333333
N11=302.175964
O11=0.4219022240948613
P11=0.4226230904617871
Q11=1527.983147437936
R11=13751.84832694142
S11=0.1465178146974578
T11=0.1469133191034391
G12=0.129484
H12=0.258968
I12=0.002964249956021043
J12=0.001978121185893458
M12=0.5804443333333333
N12=1.741333
O12=0.002431273010151717
P12=0.002435427107574628
Q12=0.07515825405733333
R12=0.450949524344
S12=7.206900913417378E-06
T12=4.817569958192599E-06
G13=0.129484
H13=0.258968
I13=0.002964249956021043
J13=0.001978121185893458
O13=0.0004752041289926495
P13=0.00047601606752829
Q13=0.01469004612266667
R13=0.08814027673599999
S13=1.408623818467479E-06
T13=9.416174680034016E-07
G14=0.129484
H14=0.258968
I14=0.002964249956021043
J14=0.001978121185893458
M14=136.1000366666667
N14=408.30011
O14=0.5700742118164518
P14=0.5710482463260632
Q14=17.62277714774667
R14=105.73666288648
S14=0.001689842457305648
T14=0.001129602634224892
G15=0.129484
H15=0.258968
I15=0.002964249956021043
J15=0.001978121185893458
M15=1.221658
N15=2.443316
O15=0.005117086949542552
P15=0.003417220037046797
Q15=0.158185164472
R15=0.6327406578879999
S15=1.516832476513736E-05
T15=6.759675352141898E-06
G16=0.129484
H16=0.258968
I16=0.002964249956021043
J16=0.001978121185893458
M16=100.7253213333333
N16=302.175964
O16=0.4219022240948613
P16=0.4226230904617871
Q16=13.04231750752533
R16=78.253905045152
S16=0.001250623649218373
T16=0.0008359996888902286
E17=1
F17=0.3333333333333333
G17=0.287247
H17=0.861741
I17=0.006575885106400611
J17=0.006582389055223097
M17=0.5804443333333333
N17=1.741333
O17=0.002431273010151717
P17=0.002435427107574628
Q17=0.166730893417
R17=1.500578040753
S17=1.598777197705046E-05
T17=1.603092873769288E-05
E18=1
F18=0.3333333333333333
G18=0.287247
H18=0.861741
I18=0.006575885106400611
J18=0.006582389055223097
O18=0.0004752041289926495
P18=0.00047601606752829
Q18=0.032588363648
R18=0.293295272832
S18=3.124887754342838E-06
T18=3.133322953008555E-06
E19=1
F19=0.3333333333333333
G19=0.287247
H19=0.861741
I19=0.006575885106400611
J19=0.006582389055223097
M19=136.1000366666667
N19=408.30011
O19=0.5700742118164518
P19=0.5710482463260632
Q19=39.09432723239
R19=351.84894509151
S19=0.003748742519026872
T19=0.003758861726621022
E20=1
F20=0.3333333333333333
G20=0.287247
H20=0.861741
I20=0.006575885106400611
J20=0.006582389055223097
M20=1.221658
N20=2.443316
O20=0.005117086949542552
P20=0.003417220037046797
Q20=0.3509175955259999
R20=2.105505573156
S20=3.36493758596538E-05
T20=2.24934717711459E-05
E21=1
F21=0.3333333333333333
G21=0.287247
H21=0.861741
I21=0.006575885106400611
J21=0.006582389055223097
M21=100.7253213333333
N21=302.175964
O21=0.4219022240948613
P21=0.4226230904617871
Q21=2829.912518795432
R21=25469.21266915889
S21=0.2713594051964024
T21=0.2720919020643175
"@

$updates -split "`n" | ForEach-Object {
    $line = $_.Trim()
    if ($line.Length -eq 0) { return }
    $parts = $line -split "=", 2
    $cellRef = $parts[0]
    $value = [double]$parts[1]
    $ws.Range($cellRef).Value = $value
}
